$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.05290172043414
$ws.Range("D2").Value = 1.059754761767796
$ws.Range("E2").Value = 1.059493850083259
$ws.Range("F2").Value = 1.069844243357623
$ws.Range("I2").Value = 1.037582526245922
$ws.Range("J2").Value = 1.057921883996734
$ws.Range("K2").Value = 1.062483402556942
$ws.Range("L2").Value = 1.062223202176085
$ws.Range("M2").Value = 1.072545665796736
$ws.Range("N2").Value = 1.023040611526442

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05407430179326
$ws.Range("D3").Value = 1.060817489896362
$ws.Range("E3").Value = 1.060545943288722
$ws.Range("F3").Value = 1.071001740108345
$ws.Range("I3").Value = 1.037749402478622
$ws.Range("J3").Value = 1.058744050205429
$ws.Range("K3").Value = 1.063359794103909
$ws.Range("L3").Value = 1.063088934069524
$ws.Range("M3").Value = 1.073518565291156
$ws.Range("N3").Value = 1.023319477699128

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.054832996574571
$ws.Range("D4").Value = 1.06150541694946
$ws.Range("E4").Value = 1.061227041852626
$ws.Range("F4").Value = 1.07175122402274
$ws.Range("I4").Value = 1.037855794540809
$ws.Range("J4").Value = 1.059275452575647
$ws.Range("K4").Value = 1.063926547301242
$ws.Range("L4").Value = 1.063648840480272
$ws.Range("M4").Value = 1.074148017199642
$ws.Range("N4").Value = 1.02349958467486

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.055151942440731
$ws.Range("D5").Value = 1.061794686698413
$ws.Range("E5").Value = 1.061513453559917
$ws.Range("F5").Value = 1.072066428701721
$ws.Range("I5").Value = 1.037900141502708
$ws.Range("J5").Value = 1.059498712263156
$ws.Range("K5").Value = 1.064164731528003
$ws.Range("L5").Value = 1.063884158440137
$ws.Range("M5").Value = 1.074412620244199
$ws.Range("N5").Value = 1.023575220743116

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.055205494333502
$ws.Range("D6").Value = 1.061843260196981
$ws.Range("E6").Value = 1.061561547921738
$ws.Range("F6").Value = 1.072119360101252
$ws.Range("I6").Value = 1.037907565251543
$ws.Range("J6").Value = 1.059536190238577
$ws.Range("K6").Value = 1.064204719096356
$ws.Range("L6").Value = 1.063923665459913
$ws.Range("M6").Value = 1.074457047184372
$ws.Range("N6").Value = 1.023587915628495

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.054837258380533
$ws.Range("D7").Value = 1.061509281929298
$ws.Range("E7").Value = 1.061230868593491
$ws.Range("F7").Value = 1.071755435324526
$ws.Range("I7").Value = 1.03785638860136
$ws.Range("J7").Value = 1.059278436339509
$ws.Range("K7").Value = 1.063929730239789
$ws.Range("L7").Value = 1.063651985072014
$ws.Range("M7").Value = 1.074151552910813
$ws.Range("N7").Value = 1.023500595645564

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053298009493053
$ws.Range("D8").Value = 1.060113859844302
$ws.Range("E8").Value = 1.059849343037764
$ws.Range("F8").Value = 1.070235320295859
$ws.Range("I8").Value = 1.037639251608605
$ws.Range("J8").Value = 1.058199862161057
$ws.Range("K8").Value = 1.062779651882739
$ws.Range("L8").Value = 1.062515838507231
$ws.Range("M8").Value = 1.072874478361892
$ws.Range("N8").Value = 1.023134925678807

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050585273396878
$ws.Range("D9").Value = 1.057657005695527
$ws.Range("E9").Value = 1.057417384061224
$ws.Range("F9").Value = 1.067560537812254
$ws.Range("I9").Value = 1.037244469693185
$ws.Range("J9").Value = 1.056294709159581
$ws.Range("K9").Value = 1.060750517148905
$ws.Range("L9").Value = 1.060511642775723
$ws.Range("M9").Value = 1.070623485284153
$ws.Range("N9").Value = 1.022487977385821

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048776455858164
$ws.Range("D10").Value = 1.056020459922997
$ws.Range("E10").Value = 1.055797711535535
$ws.Range("F10").Value = 1.065779907369253
$ws.Range("I10").Value = 1.03697311045669
$ws.Range("J10").Value = 1.055021509764531
$ws.Range("K10").Value = 1.059396020361056
$ws.Range("L10").Value = 1.059174035143658
$ws.Range("M10").Value = 1.06912237449905
$ws.Range("N10").Value = 1.022054932872313

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047993124093168
$ws.Range("D11").Value = 1.055312132530879
$ws.Range("E11").Value = 1.055096756546671
$ws.Range("F11").Value = 1.065009472980102
$ws.Range("I11").Value = 1.036853671147499
$ws.Range("J11").Value = 1.054469458994204
$ws.Range("K11").Value = 1.058809088891074
$ws.Range("L11").Value = 1.058594479591104
$ws.Range("M11").Value = 1.068472263221354
$ws.Range("N11").Value = 1.021867004640028

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.047702142988619
$ws.Range("D12").Value = 1.055049073732651
$ws.Range("E12").Value = 1.054836446176439
$ws.Range("F12").Value = 1.064723387200996
$ws.Range("I12").Value = 1.036809014726765
$ws.Range("J12").Value = 1.05426428984406
$ws.Range("K12").Value = 1.058591011821798
$ws.Range("L12").Value = 1.058379151729898
$ws.Range("M12").Value = 1.06823076412355
$ws.Range("N12").Value = 1.021797136922967

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.047764560208206
$ws.Range("D13").Value = 1.05510549868861
$ws.Range("E13").Value = 1.054892281132478
$ws.Range("F13").Value = 1.064784749609501
$ws.Range("I13").Value = 1.036818606865072
$ws.Range("J13").Value = 1.054308304443782
$ws.Range("K13").Value = 1.058637793027022
$ws.Range("L13").Value = 1.058425342791935
$ws.Range("M13").Value = 1.068282567363408
$ws.Range("N13").Value = 1.02181212662972

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047969071855931
$ws.Range("D14").Value = 1.055290387083738
$ws.Range("E14").Value = 1.055075238081398
$ws.Range("F14").Value = 1.064985823252949
$ws.Range("I14").Value = 1.036849985777351
$ws.Range("J14").Value = 1.054452501955608
$ws.Range("K14").Value = 1.058791063893979
$ws.Range("L14").Value = 1.058576681651179
$ws.Range("M14").Value = 1.068452301218831
$ws.Range("N14").Value = 1.021861230634737

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04809507595928
$ws.Range("D15").Value = 1.055404308938742
$ws.Range("E15").Value = 1.055187971233714
$ws.Range("F15").Value = 1.065109722977766
$ws.Range("I15").Value = 1.036869280753967
$ws.Range("J15").Value = 1.054541331819157
$ws.Range("K15").Value = 1.058885490554951
$ws.Range("L15").Value = 1.058669919193355
$ws.Range("M15").Value = 1.068556877332852
$ws.Range("N15").Value = 1.021891476906519

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048828440733187
$ws.Range("D16").Value = 1.056067475703037
$ws.Range("E16").Value = 1.055844239431714
$ws.Range("F16").Value = 1.065831050928857
$ws.Range("I16").Value = 1.036980996395453
$ws.Range("J16").Value = 1.055058131758858
$ws.Range("K16").Value = 1.059434964049122
$ws.Range("L16").Value = 1.059212490641982
$ws.Range("M16").Value = 1.069165517644122
$ws.Range("N16").Value = 1.022067396250465

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.049288432751671
$ws.Range("D17").Value = 1.056483544588573
$ws.Range("E17").Value = 1.056255998779228
$ws.Range("F17").Value = 1.066283678358002
$ws.Range("I17").Value = 1.037050553524411
$ws.Range("J17").Value = 1.055382106226251
$ws.Range("K17").Value = 1.059779519883125
$ws.Range("L17").Value = 1.059552733839169
$ws.Range("M17").Value = 1.069547268879121
$ws.Range("N17").Value = 1.022177634039672

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049556728589396
$ws.Range("D18").Value = 1.05672626026841
$ws.Range("E18").Value = 1.056496207160327
$ws.Range("F18").Value = 1.066547745295311
$ws.Range("I18").Value = 1.037090937940926
$ws.Range("J18").Value = 1.055571003020064
$ws.Range("K18").Value = 1.059980452469854
$ws.Range("L18").Value = 1.059751156995318
$ws.Range("M18").Value = 1.069769926348097
$ws.Range("N18").Value = 1.02224189364235

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.049648208933189
$ws.Range("D19").Value = 1.056809025108241
$ws.Range("E19").Value = 1.056578118157131
$ws.Range("F19").Value = 1.066637794993078
$ws.Range("I19").Value = 1.037104676252372
$ws.Range("J19").Value = 1.055635399748071
$ws.Range("K19").Value = 1.060048958355195
$ws.Range("L19").Value = 1.059818808256266
$ws.Range("M19").Value = 1.069845844851251
$ws.Range("N19").Value = 1.022263797685967

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.049239080977897
$ws.Range("D20").Value = 1.056438901263697
$ws.Range("E20").Value = 1.056211817166554
$ws.Range("F20").Value = 1.066235109826916
$ws.Range("I20").Value = 1.037043110053517
$ws.Range("J20").Value = 1.055347354287307
$ws.Range("K20").Value = 1.059742556548689
$ws.Range("L20").Value = 1.059516232583088
$ws.Range("M20").Value = 1.069506311786089
$ws.Range("N20").Value = 1.022165810730786

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047908848749139
$ws.Range("D21").Value = 1.055235940808273
$ws.Range("E21").Value = 1.055021360292574
$ws.Range("F21").Value = 1.064926609658031
$ws.Range("I21").Value = 1.036840753512569
$ws.Range("J21").Value = 1.054410042504862
$ws.Range("K21").Value = 1.058745931232203
$ws.Range("L21").Value = 1.058532117656039
$ws.Range("M21").Value = 1.068402319354427
$ws.Range("N21").Value = 1.021846772463971

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.047072380949392
$ws.Range("D22").Value = 1.05447985491944
$ws.Range("E22").Value = 1.054273193719156
$ws.Range("F22").Value = 1.064104411563677
$ws.Range("I22").Value = 1.03671183828349
$ws.Range("J22").Value = 1.053820063761327
$ws.Range("K22").Value = 1.058118939317592
$ws.Range("L22").Value = 1.057913046289361
$ws.Range("M22").Value = 1.067708086226026
$ws.Range("N22").Value = 1.02164581700794

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.047515818150396
$ws.Range("D23").Value = 1.054880645650562
$ws.Range("E23").Value = 1.054669780775225
$ws.Range("F23").Value = 1.064540226467072
$ws.Range("I23").Value = 1.036780338476785
$ws.Range("J23").Value = 1.054132884840133
$ws.Range("K23").Value = 1.058451355273064
$ws.Range("L23").Value = 1.0582412582114
$ws.Range("M23").Value = 1.068076122945945
$ws.Range("N23").Value = 1.021752381793797

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.049261380949738
$ws.Range("D24").Value = 1.056459073568389
$ws.Range("E24").Value = 1.056231780823669
$ws.Range("F24").Value = 1.066257055679021
$ws.Range("I24").Value = 1.03704647401559
$ws.Range("J24").Value = 1.055363057414932
$ws.Range("K24").Value = 1.059759258814872
$ws.Range("L24").Value = 1.059532726037244
$ws.Range("M24").Value = 1.069524818568907
$ws.Range("N24").Value = 1.022171153299385

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051286633217119
$ws.Range("D25").Value = 1.058291920245181
$ws.Range("E25").Value = 1.058045814054761
$ws.Range("F25").Value = 1.068251579883076
$ws.Range("I25").Value = 1.037347970785646
$ws.Range("J25").Value = 1.056787780864965
$ws.Range("K25").Value = 1.061275401940451
$ws.Range("L25").Value = 1.061030033952938
$ws.Range("M25").Value = 1.071205497922539
$ws.Range("N25").Value = 1.022655536444017
